$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 118.6930263333333
$ws.Range("H2").Value = 356.079079
$ws.Range("I2").Value = 0.2696481350657977
$ws.Range("J2").Value = 0.2696481350657977
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 63.55492266666666
$ws.Range("N2").Value = 190.664768
$ws.Range("O2").Value = 0.9936031556622397
$ws.Range("P2").Value = 0.9936031556622397
$ws.Range("Q2").Value = 7543.526109687629
$ws.Range("R2").Value = 67891.73498718865
$ws.Range("S2").Value = 0.2679232379198144
$ws.Range("T2").Value = 0.2679232379198144
$ws.Range("G3").Value = 118.6930263333333
$ws.Range("H3").Value = 356.079079
$ws.Range("I3").Value = 0.2696481350657977
$ws.Range("J3").Value = 0.2696481350657977
$ws.Range("O3").Value = 0.000612609346703606
$ws.Range("P3").Value = 0.000612609346703606
$ws.Range("Q3").Value = 4.650986236871667
$ws.Range("R3").Value = 41.858876131845
$ws.Range("S3").Value = 0.000165188967862504
$ws.Range("T3").Value = 0.000165188967862504
$ws.Range("G4").Value = 118.6930263333333
$ws.Range("H4").Value = 356.079079
$ws.Range("I4").Value = 0.2696481350657977
$ws.Range("J4").Value = 0.2696481350657977
$ws.Range("M4").Value = 0.3699833333333333
$ws.Range("N4").Value = 1.10995
$ws.Range("O4").Value = 0.005784234991056675
$ws.Range("P4").Value = 0.005784234991056675
$ws.Range("Q4").Value = 43.91444152622778
$ws.Range("R4").Value = 395.22997373605
$ws.Range("S4").Value = 0.001559708178120763
$ws.Range("T4").Value = 0.001559708178120763
$ws.Range("I5").Value = 0.3482063679522526
$ws.Range("J5").Value = 0.3482063679522526
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 63.55492266666666
$ws.Range("N5").Value = 190.664768
$ws.Range("O5").Value = 0.9936031556622397
$ws.Range("P5").Value = 0.9936031556622397
$ws.Range("Q5").Value = 9741.227498444836
$ws.Range("R5").Value = 87671.04748600352
$ws.Range("S5").Value = 0.3459789460190452
$ws.Range("T5").Value = 0.3459789460190452
$ws.Range("I6").Value = 0.3482063679522526
$ws.Range("J6").Value = 0.3482063679522526
$ws.Range("O6").Value = 0.000612609346703606
$ws.Range("P6").Value = 0.000612609346703606
$ws.Range("Q6").Value = 6.005986373841668
$ws.Range("S6").Value = 0.0002133144755892649
$ws.Range("T6").Value = 0.0002133144755892649
$ws.Range("I7").Value = 0.3482063679522526
$ws.Range("J7").Value = 0.3482063679522526
$ws.Range("M7").Value = 0.3699833333333333
$ws.Range("N7").Value = 1.10995
$ws.Range("O7").Value = 0.005784234991056675
$ws.Range("P7").Value = 0.005784234991056675
$ws.Range("Q7").Value = 56.70830314019445
$ws.Range("R7").Value = 510.37472826175
$ws.Range("S7").Value = 0.002014107457618175
$ws.Range("T7").Value = 0.002014107457618175
$ws.Range("G8").Value = 116.0670876666667
$ws.Range("H8").Value = 348.201263
$ws.Range("I8").Value = 0.2636824984472209
$ws.Range("J8").Value = 0.2636824984472209
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 63.55492266666666
$ws.Range("N8").Value = 190.664768
$ws.Range("O8").Value = 0.9936031556622397
$ws.Range("P8").Value = 0.9936031556622397
$ws.Range("Q8").Value = 7376.634780800221
$ws.Range("R8").Value = 66389.71302720199
$ws.Range("S8").Value = 0.2619957625500623
$ws.Range("T8").Value = 0.2619957625500623
$ws.Range("G9").Value = 116.0670876666667
$ws.Range("H9").Value = 348.201263
$ws.Range("I9").Value = 0.2636824984472209
$ws.Range("J9").Value = 0.2636824984472209
$ws.Range("O9").Value = 0.000612609346703606
$ws.Range("P9").Value = 0.000612609346703606
$ws.Range("Q9").Value = 4.548088830218334
$ws.Range("R9").Value = 40.93279947196501
$ws.Range("S9").Value = 0.0001615343631109266
$ws.Range("T9").Value = 0.0001615343631109266
$ws.Range("G10").Value = 116.0670876666667
$ws.Range("H10").Value = 348.201263
$ws.Range("I10").Value = 0.2636824984472209
$ws.Range("J10").Value = 0.2636824984472209
$ws.Range("M10").Value = 0.3699833333333333
$ws.Range("N10").Value = 1.10995
$ws.Range("O10").Value = 0.005784234991056675
$ws.Range("P10").Value = 0.005784234991056675
$ws.Range("Q10").Value = 42.94288798520556
$ws.Range("R10").Value = 386.48599186685
$ws.Range("S10").Value = 0.001525201534047663
$ws.Range("T10").Value = 0.001525201534047663
$ws.Range("G11").Value = 52.14473966666667
$ws.Range("H11").Value = 156.434219
$ws.Range("I11").Value = 0.1184629985347288
$ws.Range("J11").Value = 0.1184629985347288
$ws.Range("K11").Value = 3
$ws.Range("L11").Value = 1
$ws.Range("M11").Value = 63.55492266666666
$ws.Range("N11").Value = 190.664768
$ws.Range("O11").Value = 0.9936031556622397
$ws.Range("P11").Value = 0.9936031556622397
$ws.Range("Q11").Value = 3314.054896988466
$ws.Range("R11").Value = 29826.49407289619
$ws.Range("S11").Value = 0.1177052091733178
$ws.Range("T11").Value = 0.1177052091733178
$ws.Range("G12").Value = 52.14473966666667
$ws.Range("H12").Value = 156.434219
$ws.Range("I12").Value = 0.1184629985347288
$ws.Range("J12").Value = 0.1184629985347288
$ws.Range("O12").Value = 0.000612609346703606
$ws.Range("P12").Value = 0.000612609346703606
$ws.Range("Q12").Value = 2.043291623838334
$ws.Range("R12").Value = 18.389624614545
$ws.Range("S12").Value = 0.00007257154014091044
$ws.Range("T12").Value = 0.00007257154014091044
$ws.Range("G13").Value = 52.14473966666667
$ws.Range("H13").Value = 156.434219
$ws.Range("I13").Value = 0.1184629985347288
$ws.Range("J13").Value = 0.1184629985347288
$ws.Range("M13").Value = 0.3699833333333333
$ws.Range("N13").Value = 1.10995
$ws.Range("O13").Value = 0.005784234991056675
$ws.Range("P13").Value = 0.005784234991056675
$ws.Range("Q13").Value = 19.29268459767222
$ws.Range("R13").Value = 173.63416137905
$ws.Range("S13").Value = 0.0006852178212700739
$ws.Range("T13").Value = 0.000685217821270074
